# Auto-generated edit script: updates Leve-profit market-data cells
# per the scheduled-runner price refresh (H..N columns) across all 8 sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1167.8077
$ws.Range("I33").Value = 809.58826
$ws.Range("K33").Value = 809.58826
$ws.Range("M33").Value = -580.58826
$ws.Range("H40").Value = 1363.6364
$ws.Range("I40").Value = 1300
$ws.Range("J40").Value = 1440
$ws.Range("K40").Value = 1300
$ws.Range("L40").Value = 1440
$ws.Range("M40").Value = -1125
$ws.Range("N40").Value = -1790
$ws.Range("H137").Value = 1126.3667
$ws.Range("I137").Value = 1062.3684
$ws.Range("J137").Value = 1236.909
$ws.Range("K137").Value = 3187.1052
$ws.Range("L137").Value = 3710.727
$ws.Range("M137").Value = -637.1052
$ws.Range("N137").Value = -8810.727000000001
$ws.Range("H138").Value = 2533.1724
$ws.Range("I138").Value = 2854.0715
$ws.Range("J138").Value = 2431.068
$ws.Range("K138").Value = 8562.2145
$ws.Range("L138").Value = 7293.204000000001
$ws.Range("M138").Value = -3422.2145
$ws.Range("N138").Value = -17573.204

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1073.8276
$ws.Range("I2").Value = 1064
$ws.Range("J2").Value = 1104.7142
$ws.Range("K2").Value = 1064
$ws.Range("L2").Value = 1104.7142
$ws.Range("M2").Value = -951
$ws.Range("N2").Value = -1330.7142
$ws.Range("H4").Value = 17.5
$ws.Range("I4").Value = 18.333334
$ws.Range("J4").Value = 15
$ws.Range("K4").Value = 18.333334
$ws.Range("L4").Value = 15
$ws.Range("M4").Value = 97.66666599999999
$ws.Range("N4").Value = -247
$ws.Range("H5").Value = 353.33334
$ws.Range("I5").Value = 200
$ws.Range("J5").Value = 430
$ws.Range("K5").Value = 200
$ws.Range("L5").Value = 430
$ws.Range("M5").Value = -88
$ws.Range("N5").Value = -654
$ws.Range("H32").Value = 426260.47
$ws.Range("I32").Value = 454492
$ws.Range("J32").Value = 120419
$ws.Range("K32").Value = 454492
$ws.Range("L32").Value = 120419
$ws.Range("M32").Value = -454205
$ws.Range("N32").Value = -120993
$ws.Range("H45").Value = 3575.125
$ws.Range("I45").Value = 2866
$ws.Range("K45").Value = 2866
$ws.Range("M45").Value = -2489
$ws.Range("H74").Value = 784.4194
$ws.Range("I74").Value = 468
$ws.Range("J74").Value = 1222.5385
$ws.Range("K74").Value = 468
$ws.Range("L74").Value = 1222.5385
$ws.Range("M74").Value = 406
$ws.Range("N74").Value = -2970.5385
$ws.Range("H77").Value = 784.4194
$ws.Range("I77").Value = 468
$ws.Range("J77").Value = 1222.5385
$ws.Range("K77").Value = 2340
$ws.Range("L77").Value = 6112.692500000001
$ws.Range("M77").Value = 2028
$ws.Range("N77").Value = -14848.6925
$ws.Range("H116").Value = 1073.8276
$ws.Range("I116").Value = 1064
$ws.Range("J116").Value = 1104.7142
$ws.Range("K116").Value = 1064
$ws.Range("L116").Value = 1104.7142
$ws.Range("M116").Value = 1230
$ws.Range("N116").Value = -5692.7142
$ws.Range("H134").Value = 78429
$ws.Range("J134").Value = 78429
$ws.Range("L134").Value = 78429
$ws.Range("N134").Value = -88569

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1073.8276
$ws.Range("I3").Value = 1064
$ws.Range("J3").Value = 1104.7142
$ws.Range("K3").Value = 1064
$ws.Range("L3").Value = 1104.7142
$ws.Range("M3").Value = -950
$ws.Range("N3").Value = -1332.7142
$ws.Range("H4").Value = 353.33334
$ws.Range("I4").Value = 200
$ws.Range("J4").Value = 430
$ws.Range("K4").Value = 200
$ws.Range("L4").Value = 430
$ws.Range("M4").Value = -85
$ws.Range("N4").Value = -660
$ws.Range("H64").Value = 554.5
$ws.Range("I64").Value = 573
$ws.Range("J64").Value = 544.2222
$ws.Range("K64").Value = 573
$ws.Range("L64").Value = 544.2222
$ws.Range("M64").Value = -348
$ws.Range("N64").Value = -994.2222
$ws.Range("H67").Value = 554.5
$ws.Range("I67").Value = 573
$ws.Range("J67").Value = 544.2222
$ws.Range("K67").Value = 573
$ws.Range("L67").Value = 544.2222
$ws.Range("M67").Value = 207
$ws.Range("N67").Value = -2104.2222
$ws.Range("H86").Value = 4142.2856
$ws.Range("I86").Value = 2998.6667
$ws.Range("K86").Value = 2998.6667
$ws.Range("M86").Value = -1875.6667
$ws.Range("H89").Value = 4142.2856
$ws.Range("I89").Value = 2998.6667
$ws.Range("K89").Value = 14993.3335
$ws.Range("M89").Value = -9377.333500000001
$ws.Range("H107").Value = 923.0625
$ws.Range("I107").Value = 685.55554
$ws.Range("J107").Value = 1228.4286
$ws.Range("K107").Value = 685.55554
$ws.Range("L107").Value = 1228.4286
$ws.Range("M107").Value = 1234.44446
$ws.Range("N107").Value = -5068.4286
$ws.Range("H134").Value = 3023.375
$ws.Range("I134").Value = 2984.7368
$ws.Range("K134").Value = 8954.2104
$ws.Range("M134").Value = -6419.2104

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2825.0579
$ws.Range("I31").Value = 1079.174
$ws.Range("J31").Value = 3698
$ws.Range("K31").Value = 1079.174
$ws.Range("L31").Value = 3698
$ws.Range("M31").Value = -784.174
$ws.Range("N31").Value = -4288
$ws.Range("H34").Value = 2825.0579
$ws.Range("I34").Value = 1079.174
$ws.Range("J34").Value = 3698
$ws.Range("K34").Value = 1079.174
$ws.Range("L34").Value = 3698
$ws.Range("M34").Value = -877.174
$ws.Range("N34").Value = -4102
$ws.Range("H99").Value = 1963.1578
$ws.Range("I99").Value = 1500
$ws.Range("J99").Value = 1988.8889
$ws.Range("K99").Value = 1500
$ws.Range("L99").Value = 1988.8889
$ws.Range("M99").Value = -2
$ws.Range("N99").Value = -4984.8889
$ws.Range("H107").Value = 2315321.8
$ws.Range("I107").Value = 4808072.5
$ws.Range("J107").Value = 624.7143
$ws.Range("K107").Value = 4808072.5
$ws.Range("L107").Value = 624.7143
$ws.Range("M107").Value = -4806152.5
$ws.Range("N107").Value = -4464.7143
$ws.Range("H126").Value = 1963.1578
$ws.Range("I126").Value = 1500
$ws.Range("J126").Value = 1988.8889
$ws.Range("K126").Value = 4500
$ws.Range("L126").Value = 5966.6667
$ws.Range("M126").Value = -2030
$ws.Range("N126").Value = -10906.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 111111220
$ws.Range("I23").Value = 67.5
$ws.Range("J23").Value = 200000130
$ws.Range("K23").Value = 202.5
$ws.Range("L23").Value = 600000390
$ws.Range("M23").Value = 32.5
$ws.Range("N23").Value = -600000860
$ws.Range("H101").Value = 3029
$ws.Range("J101").Value = 3029
$ws.Range("L101").Value = 9087
$ws.Range("N101").Value = -13955
$ws.Range("H107").Value = 1623.591
$ws.Range("J107").Value = 2712.1943
$ws.Range("L107").Value = 8136.5829
$ws.Range("N107").Value = -11976.5829
$ws.Range("H113").Value = 1013.1818
$ws.Range("I113").Value = 557.3333
$ws.Range("J113").Value = 1560.2
$ws.Range("K113").Value = 1671.9999
$ws.Range("L113").Value = 4680.6
$ws.Range("M113").Value = 498.0001
$ws.Range("N113").Value = -9020.6
$ws.Range("H122").Value = 6435.2354
$ws.Range("I122").Value = 409.36365
$ws.Range("K122").Value = 3684.27285
$ws.Range("M122").Value = -1234.27285
$ws.Range("H129").Value = 1162.2667
$ws.Range("I129").Value = 538.9
$ws.Range("J129").Value = 1473.95
$ws.Range("K129").Value = 1616.7
$ws.Range("L129").Value = 4421.85
$ws.Range("M129").Value = 3383.3
$ws.Range("N129").Value = -14421.85
$ws.Range("H138").Value = 1725.2
$ws.Range("I138").Value = 1131.7858
$ws.Range("J138").Value = 10033
$ws.Range("K138").Value = 3395.3574
$ws.Range("L138").Value = 30099
$ws.Range("M138").Value = 1744.6426
$ws.Range("N138").Value = -40379

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1661.0834
$ws.Range("I97").Value = 1792
$ws.Range("J97").Value = 1477.8
$ws.Range("K97").Value = 1792
$ws.Range("L97").Value = 1477.8
$ws.Range("M97").Value = -1296
$ws.Range("N97").Value = -2469.8
$ws.Range("H102").Value = 2328.8572
$ws.Range("I102").Value = 2262.4
$ws.Range("K102").Value = 2262.4
$ws.Range("M102").Value = -640.4000000000001
$ws.Range("H132").Value = 2491.2173
$ws.Range("I132").Value = 2287.3333
$ws.Range("J132").Value = 2713.6365
$ws.Range("K132").Value = 6861.999899999999
$ws.Range("L132").Value = 8140.9095
$ws.Range("M132").Value = -4331.999899999999
$ws.Range("N132").Value = -13200.9095

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5507.143
$ws.Range("I122").Value = 2700
$ws.Range("J122").Value = 5975
$ws.Range("K122").Value = 8100
$ws.Range("L122").Value = 17925
$ws.Range("M122").Value = -5650
$ws.Range("N122").Value = -22825
$ws.Range("H132").Value = 3487.4333
$ws.Range("I132").Value = 3480.9285
$ws.Range("J132").Value = 3493.125
$ws.Range("K132").Value = 10442.7855
$ws.Range("L132").Value = 10479.375
$ws.Range("M132").Value = -7912.7855
$ws.Range("N132").Value = -15539.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 43809.832
$ws.Range("J46").Value = 43809.832
$ws.Range("L46").Value = 43809.832
$ws.Range("N46").Value = -44271.832
$ws.Range("H134").Value = 43809.832
$ws.Range("J134").Value = 43809.832
$ws.Range("L134").Value = 131429.496
$ws.Range("N134").Value = -136499.496
